$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final Key/Value/VoiceDuration table data for rows 2-102 (row 1 header is unchanged).
# Columns: RowNumber, A (Key), B (Value), C (VoiceDuration)
$data = @(
    @(2, "welcome", "Welcome!", 1.5),
    @(3, "title", "Project Sunflower", $null),
    @(4, "launch", "LAUNCH", $null),
    @(5, "mismatch", "MISMATCH", $null),
    @(6, "options", "OPTIONS", $null),
    @(7, "music", "MUSIC", $null),
    @(8, "sound", "SOUND", $null),
    @(9, "on", "ON", $null),
    @(10, "off", "OFF", $null),
    @(11, "close", "CLOSE", $null),
    @(12, "okay", "OKAY", $null),
    @(13, "reveal", "REVEAL", $null),
    @(14, "back", "BACK", $null),
    @(15, "levelMatchTitle", "MATCH CLIMATE", 1.5),
    @(16, "levelMatchDesc", "Find the location that matches the following climate.", 4),
    @(17, "levelMatchNotFound", "Climate does not match, try another location.", 4),
    @(18, "levelMatchFound", "Climate Match Found!", 3),
    @(19, "cycle", "CYCLE {0}", $null),
    @(20, "cycleStart", "BEGIN", $null),
    @(21, "cycleEnd", "CYCLE FINISH", $null),
    @(22, "victory", "VICTORY", $null),
    @(23, "climate", "Climate", $null),
    @(24, "climateZone", "Climate Zone", $null),
    @(25, "climateZones", "Climate Zones", $null),
    @(26, "latitudes", "Latitudes", $null),
    @(27, "latitudeEquator", "Equator (0°)", $null),
    @(28, "latitudeTropicCancer", "Tropic of Cancer (23°27')", $null),
    @(29, "latitudeTropicCapricorn", "Tropic of Capricorn (23°27')", $null),
    @(30, "latitudeArctic", "Arctic Circle (66°33')", $null),
    @(31, "latitudeAntarctic", "Antarctic Circle (66°33')", $null),
    @(32, "zone", "Zone", $null),
    @(33, "zoneTropics", "Tropics", $null),
    @(34, "zoneSubtropics", "Subtropics", $null),
    @(35, "zoneTemperate", "Temperate", $null),
    @(36, "zoneFrigid", "Frigid", $null),
    @(37, "type", "Type", $null),
    @(38, "humidity", "Humidity", $null),
    @(39, "wind", "Wind", $null),
    @(40, "weatherForecast", "Weather Forecast", $null),
    @(41, "weatherSunny", "Sunny", $null),
    @(42, "weatherSunnyDesc", "It's always sunny.", $null),
    @(43, "weatherPartlySunny", "Partly Sunny", $null),
    @(44, "weatherPartlySunnyDesc", "Sunny'ish.", $null),
    @(45, "weatherMostlyCloudy", "Mostly Cloudy", $null),
    @(46, "weatherMostlyCloudyDesc", "Clouds, clouds everywhere.", $null),
    @(47, "weatherCloudy", "Cloudy", $null),
    @(48, "weatherCloudyDesc", "Just clouds.", $null),
    @(49, "weatherClear", "Clear", $null),
    @(50, "weatherClearDesc", "No clouds allowed.", $null),
    @(51, "weatherLightRain", "Light Rain", $null),
    @(52, "weatherLightRainDesc", "Drip here and there.", $null),
    @(53, "weatherRain", "Rain", $null),
    @(54, "weatherRainDesc", "Drips.", $null),
    @(55, "weatherOvercast", "Overcast", $null),
    @(56, "weatherOvercastDesc", "No sun.", $null),
    @(57, "weatherFog", "Fog", $null),
    @(58, "weatherFogDesc", "Condensed BS", $null),
    @(59, "weatherHaze", "Haze", $null),
    @(60, "weatherHazeDesc", "Bad Omen", $null),
    @(61, "weatherSandstorm", "Sandstorm", $null),
    @(62, "weatherSandstormDesc", "danger of flying off", $null),
    @(63, "weatherLightSnow", "Light Snow", $null),
    @(64, "weatherLightSnowDesc", "a dash of snow", $null),
    @(65, "weatherSnow", "Snow", $null),
    @(66, "weatherSnowDesc", "a bunch of snow", $null),
    @(67, "weatherBlizzard", "Blizzard", $null),
    @(68, "weatherBlizzardDesc", "blizzard", $null),
    @(69, "unitAllyMallet", "Mallet Guy", $null),
    @(70, "unitAllyMalletDesc", "He has a mallet.", $null),
    @(71, "unitAllyGardener", "Gardener", $null),
    @(72, "unitAllyGardenerDesc", "He gardens stuff.", $null),
    @(73, "unitAllySpearman", "Spear Guy", $null),
    @(74, "unitAllySpearmanDesc", "He spears in air.", $null),
    @(75, "climateZonePolar", "Polar", $null),
    @(76, "climateZonePolarDesc", "Polar", $null),
    @(77, "climateZoneTemperate", "Temperate", $null),
    @(78, "climateZoneTemperateDesc", "Temperate", $null),
    @(79, "climateZoneTropical", "Tropical", $null),
    @(80, "climateZoneTropicalDesc", "Tropical", $null),
    @(81, "climateTypeArid", "Desert", 0.6),
    @(82, "climateTypeAridDesc", "Desert climates are dry areas with very little rain and humidity. Not much vegitations can grow in such a place.", 5),
    @(83, "climateTypeContinental", "Continental", 1),
    @(84, "climateTypeContinentalDesc", "Continental climates are found in inland areas. They generally have four seasons: spring, summer, autumn, and winter. Temperature varies from hot to cold as season changes.", 5),
    @(85, "climateTypeMarineWestCoast", "Oceanic", $null),
    @(86, "climateTypeMarineWestCoastDesc", "Marine West Coast", $null),
    @(87, "climateTypeMediterranean", "Mediterranean", $null),
    @(88, "climateTypeMediterraneanDesc", "Mediterranean", $null),
    @(89, "climateTypeMonsoon", "Monsoon", $null),
    @(90, "climateTypeMonsoonDesc", "Monsoon", $null),
    @(91, "climateTypePolar", "Artic", $null),
    @(92, "climateTypePolarDesc", "Artic", $null),
    @(93, "climateTypeRainforest", "Rainforest", $null),
    @(94, "climateTypeRainforestDesc", "Rainforest", $null),
    @(95, "locationGreatPlains", "Great Plains", $null),
    @(96, "locationSahara", "Sahara Desert", $null),
    @(97, "locationBritishIsles", "British Isles", $null),
    @(98, "locationPacificNorthwest", "Pacific Northwest", $null),
    @(99, "locationGreenland", "Greenland", $null),
    @(100, "tutorialLevelSelect01", "This is a satellite map of Earth. Here you will help us find a location that matches the climate we need to evaluate.", $null),
    @(101, "tutorialLevelSelect02", "You can press the image of the climate to get more information.", $null),
    @(102, "tutorialDragInstruction", "Drag the card towards the playing field.", $null)
)

foreach ($item in $data) {
    $r = $item[0]
    $a = $item[1]
    $b = $item[2]
    $c = $item[3]

    $ws.Cells.Item($r, 1).Value2 = $a
    $ws.Cells.Item($r, 2).Value2 = $b

    if ($null -ne $c) {
        $ws.Cells.Item($r, 3).Value2 = $c
    } else {
        $ws.Cells.Item($r, 3).ClearContents() | Out-Null
    }

    $ws.Cells.Item($r, 4).ClearContents() | Out-Null
}

# Row 2 keeps its MaxChars (D) value of 50; restore it since the loop above clears column D.
$ws.Cells.Item(2, 4).Value2 = 50

# Update the view state to match the edited workbook (scrolled down to the newly added rows).
$sheetView = $ws.Application.ActiveWindow
$sheetView.ScrollRow = 76
$ws.Range("B89").Select() | Out-Null

Write-Host "Edit complete"
